$d = $word.ActiveDocument

# --- Text Box 3 ("Changes for database and email ...") ---
# Fix subject/verb agreement ("has to be done" -> "must be done"),
# singular "page" instead of "pages", and add abbreviation period "i.e." (+3).
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    $tr = $shp.TextFrame.TextRange
    $txt = $tr.Text

    if ($txt -eq "Changes for database and email has to be done in every pages (i.e, 3)") {
        $tr.Text = "Changes for database and email must be done in every page (i.e., 3)"
    }
    elseif ($txt -eq "Changes are only made inside implementation inside Dependency Injection Container(DI)") {
        $tr.Text = "Changes are only made inside implementation inside Dependency Injection Container (DI)"
    }
}
